$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Student ID values in column A (rows 2-9)
$ws.Range("A2").Value = "200958"
$ws.Range("A3").Value = "200785"
$ws.Range("A4").Value = "201023"
$ws.Range("A5").Value = "201838"
$ws.Range("A6").Value = "211146"
$ws.Range("A7").Value = "201574"
$ws.Range("A8").Value = "201237"
$ws.Range("A9").Value = "211137"

# Remove the last row (row 10), which no longer exists in the updated log
$ws.Rows.Item(10).Delete()
